$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.061.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.497.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.44"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.83"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.06"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0803"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.889.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.499.50"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.831"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.917.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.96"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0933"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.18"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.65"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.76"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +10.10%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.71"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.89"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.36"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.14"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.28"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0772"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.55"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "121.72"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.09"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0304"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.000.59"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.15"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.90"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.17"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.72"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.06%  "
